$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Forecast Results" ---
$ws1 = $wb.Worksheets.Item("Forecast Results")

$data = @(
    @{ Row = 2;  D = 0; E = 85.6014294952854;  F = -56.6014294952854 }
    @{ Row = 3;  D = 0; E = 72.31260055855948; F = -48.31260055855948 }
    @{ Row = 4;  D = 0; E = 110.2528909564889; F = -79.25289095648888 }
    @{ Row = 5;  D = 0; E = 74.87162204434462; F = -45.87162204434462 }
    @{ Row = 6;  D = 0; E = 69.31789672524627; F = -41.31789672524627 }
    @{ Row = 7;  D = 0; E = 63.11636771285703; F = -38.11636771285703 }
    @{ Row = 8;  D = 0; E = 77.59542670934502; F = -47.59542670934502 }
    @{ Row = 9;  D = 0; E = 77.17464788917684; F = -51.17464788917684 }
    @{ Row = 10; D = 0; E = 79.00145279656964; F = -52.00145279656964 }
    @{ Row = 11; D = 0; E = 98.47482414888874; F = -58.47482414888874 }
    @{ Row = 12; D = 0; E = 77.98591913016988; F = -45.98591913016988 }
    @{ Row = 13; D = 0; E = 88.80592774512995; F = -65.80592774512995 }
)

foreach ($row in $data) {
    $ws1.Cells.Item($row.Row, 4).Value = $row.D
    $ws1.Cells.Item($row.Row, 5).Value = $row.E
    $ws1.Cells.Item($row.Row, 6).Value = $row.F
}

# --- Sheet 2: "Metrics" ---
$ws2 = $wb.Worksheets.Item("Metrics")

$ws2.Cells.Item(2, 1).Value = 52.54258382600514
$ws2.Cells.Item(2, 2).Value = 186.2018732181669
$ws2.Cells.Item(2, 3).Value = 53.64686381954046
